$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> (nombre_aides[C], montant_total[E]) updates
$updates = @{
    7   = @{ C = 7017;   E = 290920261 }
    64  = @{ C = 5219;   E = 20481621 }
    91  = @{ C = 151203; E = 482928044 }
    92  = @{ C = 409295; E = 1597281123 }
    93  = @{ C = 209658; E = 1310059700 }
    94  = @{ C = 94235;  E = 919215690 }
    95  = @{ C = 50805;  E = 934508381 }
    96  = @{ C = 17323;  E = 797283923 }
    97  = @{ C = 2163;   E = 214451252 }
    104 = @{ C = 135301; E = 272651977 }
    128 = @{ C = 25;     E = 509896 }
    132 = @{ C = 30292;  E = 174228945 }
    135 = @{ C = 1857;   E = 65742275 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 5).Value = $vals.E
}
